# Update "Pais" (COVID-19 countries) worksheet with the 25-Apr-2020 20:52 data refresh.
# Source data is re-sorted by total cases each refresh; a handful of countries
# leapfrogged their neighbours in the ranking (Peru overtakes India/Portugal/Ecuador,
# Sudafrica overtakes Egipto), which is why whole rows of figures shift down by one
# slot while a few other rows just get their counters bumped in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 20:52"

# --- Canada (row 16) - in-place tweak ---------------------------------
$ws.Range("E16").Value = 27085
$ws.Range("G16").Value = 160
$ws.Range("H16").Value = 2462

# --- Suiza (row 18) - in-place tweak -----------------------------------
$ws.Range("E18").Value = 6295
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 1599

# --- Peru overtakes India, Portugal & Ecuador (rows 19-22) -------------
# Row 19: was India -> now Peru, with Peru's freshly updated totals
$ws.Range("A19").Value = "Peru"
$ws.Range("B19").Value = 25331
$ws.Range("C19").Value = 3683
$ws.Range("D19").Value = 7797
$ws.Range("E19").Value = 16834
$ws.Range("F19").Value = 545
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = 700

# Row 20: was Portugal -> now India, carrying India's previous totals
$ws.Range("A20").Value = "India"
$ws.Range("B20").Value = 24942
$ws.Range("C20").Value = 495
$ws.Range("D20").Value = 5498
$ws.Range("E20").Value = 18664
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 780

# Row 21: was Ecuador -> now Portugal, carrying Portugal's previous totals
$ws.Range("A21").Value = "Portugal"
$ws.Range("B21").Value = 23392
$ws.Range("C21").Value = 595
$ws.Range("D21").Value = 1277
$ws.Range("E21").Value = 21235
$ws.Range("F21").Value = 186
$ws.Range("G21").Value = 26
$ws.Range("H21").Value = 880

# Row 22: was Peru -> now Ecuador, carrying Ecuador's previous totals
$ws.Range("A22").Value = "Ecuador"
$ws.Range("B22").Value = 22719
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 1366
$ws.Range("E22").Value = 20777
$ws.Range("F22").Value = 127
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 576

# --- Sudafrica overtakes Egipto (rows 53-54) ----------------------------
# Row 53: was Egipto -> now Sudafrica, with Sudafrica's freshly updated totals
$ws.Range("A53").Value = "Sudafrica"
$ws.Range("B53").Value = 4361
$ws.Range("C53").Value = 141
$ws.Range("D53").Value = 1473
$ws.Range("E53").Value = 2802
$ws.Range("F53").Value = 36
$ws.Range("G53").Value = 7
$ws.Range("H53").Value = 86

# Row 54: was Sudafrica -> now Egipto, carrying Egipto's previous totals
$ws.Range("A54").Value = "Egipto"
$ws.Range("B54").Value = 4319
$ws.Range("C54").Value = 227
$ws.Range("D54").Value = 1114
$ws.Range("E54").Value = 2898
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 13
$ws.Range("H54").Value = 307

# --- Uzbekistan (row 68) - in-place tweak --------------------------------
$ws.Range("B68").Value = 1862
$ws.Range("C68").Value = 58
$ws.Range("E68").Value = 1147

# --- Gibraltar (row 139) - in-place tweak --------------------------------
$ws.Range("B139").Value = 136
$ws.Range("C139").Value = 3
$ws.Range("E139").Value = 5

# --- Monaco (row 149) - in-place tweak -----------------------------------
$ws.Range("D149").Value = 42
$ws.Range("E149").Value = 48
$ws.Range("F149").Value = 1
